$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 2343
    $ws.Range("F3").Value = 1841
    $ws.Range("F5").Value = 1133
    $ws.Range("F6").Value = 1082
    $ws.Range("F8").Value = 5939
}

$wb.Save()
